$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
for ($r = 8; $r -le 19; $r++) {
    $ws.Range("L$r").Value = "Unbegrenzt"
}
$ws.Rows("8:19").AutoFit()
